$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The D2 cell currently holds the number 299 (price column).
# Change it to the text value "350" (stored as a string, not a number,
# like the other textual columns in this sheet).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "350"
$ws.Range("D2").Style = "Normal"
